$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("unmatched_expenses")

$ws.Range("A6").Value = "UNKNOWN DATE"
$ws.Range("B6").Value = 780
$ws.Range("C6").Value = "Amount in Words: 780 Rupees Only (parsing failed)"
